# ButecoCalc/Exemplo.xlsx — 1st delivery: populate "Plan2" with a second
# (smaller) event calculation, mirroring the layout already present on
# "Plan1", then leave Plan2 as the active/selected sheet.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Plan1")
$ws2 = $wb.Worksheets.Item("Plan2")

# ---------------------------------------------------------------------
# Style colours (BGR, as read back off Plan1 via COM .Interior.Color)
# ---------------------------------------------------------------------
$YELLOW = 65535    # FFFF00
$ORANGE = 49407    # FFC000
$NONE   = -4142    # xlColorIndexNone

$CENTER = -4108
$BOTTOM = -4107

# ---------------------------------------------------------------------
# Small helpers so every cell/range only has to state what differs from
# a plain, unformatted cell.
# ---------------------------------------------------------------------
function Set-BoldItalic($rng, [bool]$bold, [bool]$italic, $size) {
    $rng.Font.Bold = $bold
    $rng.Font.Italic = $italic
    if ($size) { $rng.Font.Size = $size }
}

function Set-Fill($rng, $color) {
    if ($color -eq $NONE) {
        $rng.Interior.ColorIndex = -4142
    } else {
        $rng.Interior.Color = $color
    }
}

# ---------------------------------------------------------------------
# Values + formulas
# ---------------------------------------------------------------------

# Row 4 — merged event title
$ws2.Range("F4").Value = "Quem"
$ws2.Range("F4:H4").Merge()

# Row 5 — column headers
$ws2.Range("C5").Value = "Qtd"
$ws2.Range("D5").Value = "Preco unitario"
$ws2.Range("E5").Value = "Total"
$ws2.Range("F5").Value = "A"
$ws2.Range("G5").Value = "B"
$ws2.Range("H5").Value = "C"
$ws2.Range("I5").Value = "Valor Individual"
$ws2.Range("J5").Value = "Valor por cabeça"

# Row 6 — cerveja
$ws2.Range("B6").Value = "cerveja"
$ws2.Range("C6").Value = 10
$ws2.Range("D6").Value = 5
$ws2.Range("E6").Formula = "=D6*C6"
$ws2.Range("F6").Value = "x"
$ws2.Range("G6").Value = "x"
$ws2.Range("H6").Value = "x"
$ws2.Range("I6").Formula = '=COUNTIF(F6:H6,"x")'
$ws2.Range("J6").Formula = "=E6/I6"

# Row 7 — tropeiro simples
$ws2.Range("B7").Value = "tropeiro simples"
$ws2.Range("C7").Value = 1
$ws2.Range("D7").Value = 18
$ws2.Range("F7").Value = "x"
$ws2.Range("G7").Value = "x"
$ws2.Range("I7").Formula = '=COUNTIF(F7:H7,"x")'
$ws2.Range("J7").Formula = "=E7/I7"

# Row 8 — caipirinha
$ws2.Range("B8").Value = "caipirinha"
$ws2.Range("C8").Value = 3
$ws2.Range("D8").Value = 8
$ws2.Range("G8").Value = "x"
$ws2.Range("H8").Value = "x"
$ws2.Range("I8").Formula = '=COUNTIF(F8:H8,"x")'
$ws2.Range("J8").Formula = "=E8/I8"

# Shared formula E7:E8 (D*C), entered on the whole block as Plan1 does
$ws2.Range("E7:E8").Formula = "=D7*C7"

# Row 9/10 — totals
$ws2.Range("D9").Value = "Total com serviço"
$ws2.Range("E9").Formula = "=SUM(E6:E8)*(1+`$C`$12/100)"
$ws2.Range("D10").Value = "Total sem serviço"
$ws2.Range("E10").Formula = "=SUM(E6:E8)"

# Row 12 — serviço % and "Valor Individual (com 10%)" breakdown
$ws2.Range("B12").Value = "Serviço"
$ws2.Range("C12").Value = 10
$ws2.Range("E12").Value = "Total Individual (com 10%)"
$ws2.Range("F12").Formula = "=(E6/I6+E7/I7)*(1+`$C`$12/100)"
$ws2.Range("G12").Formula = "=(E6/I6+E7/I7+E8/I8)*(1+`$C`$12/100)"
$ws2.Range("H12").Formula = "=(E6/I6+E8/I8)*(1+`$C`$12/100)"
$ws2.Range("I12").Formula = "=SUM(F12:H12)"

# Row 13 — "Total Individual (sem 10%)" breakdown
$ws2.Range("E13").Value = "Total Individual (sem 10%)"
$ws2.Range("F13").Formula = "=(E6/I6+E7/I7)"
$ws2.Range("G13").Formula = "=(E6/I6+E7/I7+E8/I8)"
$ws2.Range("H13").Formula = "=(E6/I6+E8/I8)"
$ws2.Range("I13").Formula = "=SUM(F13:H13)"

# ---------------------------------------------------------------------
# Formatting — mirrors the xf combinations already used on Plan1
# ---------------------------------------------------------------------

# title "Quem" — bold, size 12, centered (matches Plan1 F4:M4)
$r = $ws2.Range("F4:H4")
Set-BoldItalic $r $true $false 12
$r.HorizontalAlignment = $CENTER
$r.VerticalAlignment = $CENTER

# header row — bold, centered
$r = $ws2.Range("C5:H5")
Set-BoldItalic $r $true $false 11
$r.HorizontalAlignment = $CENTER
$r.VerticalAlignment = $CENTER

# header row, yellow highlighted cells (Total / Valor Individual / Valor por cabeça)
$r = $ws2.Range("E5")
Set-BoldItalic $r $true $false 11
$r.HorizontalAlignment = $CENTER
$r.VerticalAlignment = $CENTER
Set-Fill $r $YELLOW

$r = $ws2.Range("I5:J5")
Set-BoldItalic $r $true $false 11
$r.HorizontalAlignment = $CENTER
$r.VerticalAlignment = $CENTER
Set-Fill $r $YELLOW

# item-name labels — bold italic, default alignment
$r = $ws2.Range("B6:B8")
Set-BoldItalic $r $true $true 11

# plain numeric / mark cells — centered, no bold
$r = $ws2.Range("C6:D8")
Set-BoldItalic $r $false $false 11
$r.HorizontalAlignment = $CENTER
$r.VerticalAlignment = $CENTER

$r = $ws2.Range("F6:H8")
Set-BoldItalic $r $false $false 11
$r.HorizontalAlignment = $CENTER
$r.VerticalAlignment = $CENTER

# "Total" column (E6:E8) + totals E9:E10 — orange fill, centered
$r = $ws2.Range("E6:E8")
Set-BoldItalic $r $false $false 11
$r.HorizontalAlignment = $CENTER
$r.VerticalAlignment = $CENTER
Set-Fill $r $ORANGE

$r = $ws2.Range("E9:E10")
Set-BoldItalic $r $false $false 11
$r.HorizontalAlignment = $CENTER
$r.VerticalAlignment = $CENTER
Set-Fill $r $ORANGE

# COUNTIF column (I6:I8) — orange fill, horizontal-center only
$r = $ws2.Range("I6:I8")
Set-BoldItalic $r $false $false 11
$r.HorizontalAlignment = $CENTER
$r.VerticalAlignment = $BOTTOM
Set-Fill $r $ORANGE

# per-person value column (J6:J8) and the F12:I13 breakdown block —
# orange fill, centered, 2-decimal number format
$r = $ws2.Range("J6:J8")
Set-BoldItalic $r $false $false 11
$r.HorizontalAlignment = $CENTER
$r.VerticalAlignment = $CENTER
$r.NumberFormat = "0.00"
Set-Fill $r $ORANGE

$r = $ws2.Range("F12:I13")
Set-BoldItalic $r $false $false 11
$r.HorizontalAlignment = $CENTER
$r.VerticalAlignment = $CENTER
$r.NumberFormat = "0.00"
Set-Fill $r $ORANGE

# "Total com/sem serviço" row labels — bold italic, yellow fill, default align
$r = $ws2.Range("D9:D10")
Set-BoldItalic $r $true $true 11
Set-Fill $r $YELLOW

# "Serviço" label — bold italic, centered
$r = $ws2.Range("B12")
Set-BoldItalic $r $true $true 11
$r.HorizontalAlignment = $CENTER
$r.VerticalAlignment = $CENTER

# plain centered cells: C12, E11, J12, J13
$r = $ws2.Range("C12")
$r.HorizontalAlignment = $CENTER
$r.VerticalAlignment = $CENTER

foreach ($addr in @("E11", "J12", "J13")) {
    $r = $ws2.Range($addr)
    $r.HorizontalAlignment = $CENTER
    $r.VerticalAlignment = $CENTER
}

# "Total Individual (com/sem 10%)" labels — bold italic, centered, wrap, yellow
$r = $ws2.Range("E12:E13")
Set-BoldItalic $r $true $true 11
$r.HorizontalAlignment = $CENTER
$r.VerticalAlignment = $CENTER
$r.WrapText = $true
Set-Fill $r $YELLOW

# ---------------------------------------------------------------------
# Layout — column widths, row heights, merged title
# ---------------------------------------------------------------------
$ws2.Columns.Item(4).ColumnWidth = 16.17   # -> raw width 17   (bestFit, like Plan1 col D)
$ws2.Columns.Item(5).ColumnWidth = 11.17   # -> raw width 12
$ws2.Columns.Item(10).ColumnWidth = 14.1   # -> raw width ~15.71 (bestFit, like Plan1 col O)

$ws2.Rows.Item(4).RowHeight = 15.75
$ws2.Rows.Item(12).RowHeight = 45
$ws2.Rows.Item(13).RowHeight = 60

# ---------------------------------------------------------------------
# Selection / active sheet — Plan1 keeps a plain A1:O14 selection,
# Plan2 becomes the active tab with I13 selected.
# ---------------------------------------------------------------------
$ws1.Range("A1:O14").Select()
$ws2.Activate()
$ws2.Range("I13").Select()
